# Update "合肥-漫展信息.xlsx": remove the expired "六安lovelive only" event
# (row 2, id=81146) from the "展览" and "全部类型" sheets, shifting every
# later row up by one, and refresh the data (dates/prices/counts/links/
# covers) for the remaining events to the latest scrape.

$wb = $excel.ActiveWorkbook

# Final state (after the old row 2 is removed) for rows 2..9, columns A..I.
$rows = @(
    @{ A=1; B="2024-02-13"; C="合肥·新春AG动漫游戏盛典热血plus"; D="山西路与太原路交叉口 挥动体育"; E="2024.02.13 09:30-02.14 16:00"; F=2154; G=55; H="https://show.bilibili.com/platform/detail.html?id=80584"; I="//i1.hdslb.com/bfs/openplatform/202401/yI94srFk1704703809648.jpeg" },
    @{ A=2; B="2024-02-17"; C="合肥·2024运动新春动漫庆典（全ip）"; D="锦绣大道与清潭路交口东北角 李宁体育公园"; E="2024.02.17 09:00-02.17 17:00"; F=1639; G=65; H="https://show.bilibili.com/platform/detail.html?id=79918"; I="//i0.hdslb.com/bfs/openplatform/202312/vzuMc0sJ1702902061660.jpeg" },
    @{ A=3; B="2024-02-19"; C="合肥·安徽马娘only"; D="桐城路与庐江路交叉口西南80米 赤阑桥文玩大厦"; E="2024.02.19 09:00-02.19 17:00"; F=322; G=68; H="https://show.bilibili.com/platform/detail.html?id=78286"; I="//i1.hdslb.com/bfs/openplatform/202311/721L5pIZ1699428443216.jpeg" },
    @{ A=4; B="2024-03-02"; C="合肥·星芒1.5动漫嘉年华"; D="山西路与太原路交叉口 挥动体育"; E="2024.03.02 09:30-03.02 17:30"; F=1046; G=55; H="https://show.bilibili.com/platform/detail.html?id=81267"; I="//i0.hdslb.com/bfs/openplatform/202401/GWidiefU1706003134747.jpeg" },
    @{ A=5; B="2024-03-16"; C="合肥·CW国潮动漫游戏嘉年华"; D="南京路与庐州大道交汇处 合肥滨湖国际会展中心"; E="2024.03.16 09:30-03.17 17:00"; F=568; G=65; H="https://show.bilibili.com/platform/detail.html?id=81284"; I="//i0.hdslb.com/bfs/openplatform/202401/38B92fWF1705995243803.jpeg" },
    @{ A=6; B="2024-03-23"; C="合肥·原&铁&崩 only展"; D="金寨路与天堂窄路交叉口 梵木艺术中心"; E="2024.03.23 09:00-03.23 17:00"; F=30; G=58; H="https://show.bilibili.com/platform/detail.html?id=81574"; I="//i2.hdslb.com/bfs/openplatform/202401/0V5uyX6C1706697212904.png" },
    @{ A=7; B="2024-04-04"; C="合肥· 第二届漫画城市动漫展 -故事再次开始"; D="凤淮路与固镇路西北角 庐阳全民健身中心"; E="2024.04.04 09:00-04.05 17:00"; F=5700; G=60; H="https://show.bilibili.com/platform/detail.html?id=78898"; I="//i2.hdslb.com/bfs/openplatform/202402/QnupNcrS1707125949328.jpeg" },
    @{ A=8; B="2024-05-18"; C="合肥·梦时空SPO1动漫展"; D="阜阳路16号 银瑞林国际大酒店"; E="2024.05.18 10:00-05.18 17:00"; F=81; G=60; H="https://show.bilibili.com/platform/detail.html?id=80207"; I="//i2.hdslb.com/bfs/openplatform/202312/tQQOHYE01703574162111.jpeg" }
)

# Sheets "展览" and "全部类型" both carry this same table.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Drop the obsolete first event row; everything below shifts up.
    $ws.Rows.Item(2).Delete()

    $lastRow = 1 + $rows.Count

    # Column B holds plain "yyyy-mm-dd" text in the source data, not real
    # dates -- force text format first so Excel doesn't auto-convert the
    # assigned strings into date serials, then restore the default style
    # (the source cells carry no explicit style) once the values are set.
    $dateRange = $ws.Range("B2:B$lastRow")
    $dateRange.NumberFormat = "@"

    # Re-write the remaining 8 data rows with the refreshed values.
    for ($i = 0; $i -lt $rows.Count; $i++) {
        $r = $i + 2
        $row = $rows[$i]
        $ws.Cells.Item($r, 1).Value = $row.A
        $ws.Cells.Item($r, 2).Value = $row.B
        $ws.Cells.Item($r, 3).Value = $row.C
        $ws.Cells.Item($r, 4).Value = $row.D
        $ws.Cells.Item($r, 5).Value = $row.E
        $ws.Cells.Item($r, 6).Value = $row.F
        $ws.Cells.Item($r, 7).Value = $row.G
        $ws.Cells.Item($r, 8).Value = $row.H
        $ws.Cells.Item($r, 9).Value = $row.I
    }

    $dateRange.Style = "Normal"
}
